$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.149199999999997
$ws.Range("D10").Value = -8.167899999999989
$ws.Range("D12").Value = -8.121099999999998
$ws.Range("E13").Value = 12.44769999999999
$ws.Range("D18").Value = -8.182499999999994
$ws.Range("D25").Value = -8.328899999999997
